$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.883183028644055
$ws.Range("C2").Value = 0.2358512843325684
$ws.Range("D2").Value = 0.5864826086597219
$ws.Range("E2").Value = 0.2171666956241971
$ws.Range("G2").Value = 0.00246720104783973
$ws.Range("I2").Value = 0.7222800043408846
$ws.Range("J2").Value = 0.09812368452844211
$ws.Range("M2").Value = 0.4786273770703033
$ws.Range("N2").Value = 1.365034447659678
$ws.Range("O2").Value = 4.030090296124257
$ws.Range("B3").Value = 0.8033340464549781
$ws.Range("C3").Value = 0.2129154315258859
$ws.Range("D3").Value = 0.5826259973420918
$ws.Range("E3").Value = 0.2167951035967413
$ws.Range("G3").Value = 0.002470526933855391
$ws.Range("I3").Value = 0.7274079360384498
$ws.Range("J3").Value = 0.0985693239473342
$ws.Range("M3").Value = 0.4553933785695961
$ws.Range("N3").Value = 1.380341581171987
$ws.Range("O3").Value = 4.026415128617487
$ws.Range("B4").Value = 0.7544733867581499
$ws.Range("C4").Value = 0.1988699658759856
$ws.Range("D4").Value = 0.5805395360943351
$ws.Range("E4").Value = 0.2166678479860664
$ws.Range("G4").Value = 0.002472678883268816
$ws.Range("I4").Value = 0.7309916871335851
$ws.Range("J4").Value = 0.09888446746873214
$ws.Range("M4").Value = 0.4413221926149689
$ws.Range("N4").Value = 1.390273704016664
$ws.Range("O4").Value = 4.026647366315103
$ws.Range("B5").Value = 0.7346051610903146
$ws.Range("C5").Value = 0.1931558587390327
$ws.Range("D5").Value = 0.5797601275594246
$ws.Range("E5").Value = 0.2166413734066346
$ws.Range("G5").Value = 0.002473583526355798
$ws.Range("I5").Value = 0.7325613870678822
$ws.Range("J5").Value = 0.09902333142941799
$ws.Range("M5").Value = 0.4356372244518738
$ws.Range("N5").Value = 1.394455273696494
$ws.Range("O5").Value = 4.027367124982561
$ws.Range("B6").Value = 0.7313086747641933
$ws.Range("C6").Value = 0.1922076164812268
$ws.Range("D6").Value = 0.5796349875803628
$ws.Range("E6").Value = 0.2166385107076536
$ws.Range("G6").Value = 0.002473735417472893
$ws.Range("I6").Value = 0.7328286315246118
$ws.Range("J6").Value = 0.09904702026555157
$ws.Range("M6").Value = 0.4346962152890868
$ws.Range("N6").Value = 1.395157722533089
$ws.Range("O6").Value = 4.027524375425628
$ws.Range("B7").Value = 0.7542052620493678
$ws.Range("C7").Value = 0.198792864599568
$ws.Range("D7").Value = 0.5805287378135802
$ws.Range("E7").Value = 0.2166673881515209
$ws.Range("G7").Value = 0.00247269097135975
$ws.Range("I7").Value = 0.7310124143199808
$ws.Range("J7").Value = 0.09888629796344617
$ws.Range("M7").Value = 0.4412453238229475
$ws.Range("N7").Value = 1.390329554963468
$ws.Range("O7").Value = 4.026654543097891
$ws.Range("B8").Value = 0.855616972087148
$ws.Range("C8").Value = 0.2279353672400362
$ws.Range("D8").Value = 0.585094449511999
$ws.Range("E8").Value = 0.2170176335504088
$ws.Range("G8").Value = 0.002468325067138233
$ws.Range("I8").Value = 0.7239577101538508
$ws.Range("J8").Value = 0.0982687236438089
$ws.Range("M8").Value = 0.4705760350758155
$ws.Range("N8").Value = 1.370201624179238
$ws.Range("O8").Value = 4.028306105018487
$ws.Range("B9").Value = 1.05577988899438
$ws.Range("C9").Value = 0.2853747323190134
$ws.Range("D9").Value = 0.5962803783593529
$ws.Range("E9").Value = 0.2185049491711979
$ws.Range("G9").Value = 0.002460631194715745
$ws.Range("I9").Value = 0.7135829794006412
$ws.Range("J9").Value = 0.09738714762108813
$ws.Range("M9").Value = 0.5296309101593906
$ws.Range("N9").Value = 1.33496337111896
$ws.Range("O9").Value = 4.051330721737202
$ws.Range("B10").Value = 1.203604301136863
$ws.Range("C10").Value = 0.3277512097278077
$ws.Range("D10").Value = 0.6058600221793426
$ws.Range("E10").Value = 0.2200859490620353
$ws.Range("G10").Value = 0.002455501997094114
$ws.Range("I10").Value = 0.7080792968896219
$ws.Range("J10").Value = 0.0969404645468579
$ws.Range("M10").Value = 0.5739520808066416
$ws.Range("N10").Value = 1.311653108080243
$ws.Range("O10").Value = 4.080371132503757
$ws.Range("B11").Value = 1.271015269167378
$ws.Range("C11").Value = 0.3470674782654726
$ws.Range("D11").Value = 0.6105138487145609
$ws.Range("E11").Value = 0.2209113041731179
$ws.Range("G11").Value = 0.002453281102576841
$ws.Range("I11").Value = 0.7060376953966596
$ws.Range("J11").Value = 0.09678094476365473
$ws.Range("M11").Value = 0.5943171882322247
$ws.Range("N11").Value = 1.301608592620383
$ws.Range("O11").Value = 4.096229097407814
$ws.Range("B12").Value = 1.296565064724632
$ws.Range("C12").Value = 0.3543875503927723
$ws.Range("D12").Value = 0.6123186777483909
$ws.Range("E12").Value = 0.2212391087770165
$ws.Range("G12").Value = 0.002452456183020364
$ws.Range("I12").Value = 0.7053312055741401
$ws.Range("J12").Value = 0.09672682234773333
$ws.Range("M12").Value = 0.602057994086806
$ws.Range("N12").Value = 1.29788546793587
$ws.Range("O12").Value = 4.102615763873018
$ws.Range("B13").Value = 1.291061464087022
$ws.Range("C13").Value = 0.3528108037795334
$ws.Range("D13").Value = 0.6119280847735524
$ws.Range("E13").Value = 0.2211678315307744
$ws.Range("G13").Value = 0.002452633130091919
$ws.Range("I13").Value = 0.7054803949975934
$ws.Range("J13").Value = 0.09673819902939584
$ws.Range("M13").Value = 0.6003895878225194
$ws.Range("N13").Value = 1.298683727367376
$ws.Range("O13").Value = 4.101223294611799
$ws.Range("B14").Value = 1.273116816424476
$ws.Range("C14").Value = 0.3476695974387098
$ws.Range("D14").Value = 0.6106614810187239
$ws.Range("E14").Value = 0.220937967031773
$ws.Range("G14").Value = 0.002453212914129734
$ws.Range("I14").Value = 0.7059782358607336
$ws.Range("J14").Value = 0.09677636611421292
$ws.Range("M14").Value = 0.5949534495004087
$ws.Range("N14").Value = 1.301300674855096
$ws.Range("O14").Value = 4.096746879222366
$ws.Range("B15").Value = 1.26212813181354
$ws.Range("C15").Value = 0.3445211604737324
$ws.Range("D15").Value = 0.6098911866540675
$ws.Range("E15").Value = 0.2207991556546212
$ws.Range("G15").Value = 0.002453570140887071
$ws.Range("I15").Value = 0.7062918589950087
$ws.Range("J15").Value = 0.09680056305387907
$ws.Range("M15").Value = 0.5916274250969025
$ws.Range("N15").Value = 1.302914118342876
$ws.Range("O15").Value = 4.09405466829034
$ws.Range("B16").Value = 1.199202069553507
$ws.Range("C16").Value = 0.3264896148296543
$ws.Range("D16").Value = 0.6055618366261513
$ws.Range("E16").Value = 0.2200341458533899
$ws.Range("G16").Value = 0.002455649393053327
$ws.Range("I16").Value = 0.7082220331948506
$ws.Range("J16").Value = 0.09695176874721767
$ws.Range("M16").Value = 0.5726252415273905
$ws.Range("N16").Value = 1.312320803898281
$ws.Range("O16").Value = 4.079388126535747
$ws.Range("B17").Value = 1.160640499663828
$ws.Range("C17").Value = 0.3154376999663953
$ws.Range("D17").Value = 0.6029817054465525
$ws.Range("E17").Value = 0.21959202200593
$ws.Range("G17").Value = 0.002456953679964175
$ws.Range("I17").Value = 0.7095246013354526
$ws.Range("J17").Value = 0.09705571800525448
$ws.Range("M17").Value = 0.5610198894656477
$ws.Range("N17").Value = 1.318234834040421
$ws.Range("O17").Value = 4.071069374553844
$ws.Range("B18").Value = 1.13847648465287
$ws.Range("C18").Value = 0.3090846124942459
$ws.Range("D18").Value = 0.6015255452614952
$ws.Range("E18").Value = 0.2193477156416535
$ws.Range("G18").Value = 0.002457714455759785
$ws.Range("I18").Value = 0.7103172838365239
$ws.Range("J18").Value = 0.09711961789577828
$ws.Range("M18").Value = 0.5543639366331874
$ws.Range("N18").Value = 1.321689080486227
$ws.Range("O18").Value = 4.066533773190912
$ws.Range("B19").Value = 1.130974836059522
$ws.Range("C19").Value = 0.3069342042859091
$ws.Range("D19").Value = 0.6010373011236823
$ws.Range("E19").Value = 0.2192667138866717
$ws.Range("G19").Value = 0.002457973861763269
$ws.Range("I19").Value = 0.7105931347961345
$ws.Range("J19").Value = 0.09714195927453417
$ws.Range("M19").Value = 0.5521136396776711
$ws.Range("N19").Value = 1.322867670135409
$ws.Range("O19").Value = 4.065040856117946
$ws.Range("B20").Value = 1.164743839296591
$ws.Range("C20").Value = 0.3166138154954297
$ws.Range("D20").Value = 0.6032534813767825
$ws.Range("E20").Value = 0.2196380527816011
$ws.Range("G20").Value = 0.002456813741441329
$ws.Range("I20").Value = 0.7093814394521232
$ws.Range("J20").Value = 0.09704422692919579
$ws.Range("M20").Value = 0.562253319932978
$ws.Range("N20").Value = 1.317599825590776
$ws.Range("O20").Value = 4.071929129794967
$ws.Range("B21").Value = 1.278386985528243
$ws.Range("C21").Value = 0.3491795487416312
$ws.Range("D21").Value = 0.6110323594045042
$ws.Range("E21").Value = 0.2210050696690296
$ws.Range("G21").Value = 0.002453042181542949
$ws.Range("I21").Value = 0.7058301984519346
$ws.Range("J21").Value = 0.09676498493602281
$ws.Range("M21").Value = 0.5965493904203356
$ws.Range("N21").Value = 1.300529828208763
$ws.Range("O21").Value = 4.098051347115273
$ws.Range("B22").Value = 1.352791354319606
$ws.Range("C22").Value = 0.3704946465084618
$ws.Range("D22").Value = 0.616364162280945
$ws.Range("E22").Value = 0.221987442773468
$ws.Range("G22").Value = 0.002450670968855879
$ws.Range("I22").Value = 0.7038976547143818
$ws.Range("J22").Value = 0.09661911530166023
$ws.Range("M22").Value = 0.6191326737419018
$ws.Range("N22").Value = 1.289842938325336
$ws.Range("O22").Value = 4.117348407393649
$ws.Range("B23").Value = 1.313068589321801
$ws.Range("C23").Value = 0.3591155561785513
$ws.Range("D23").Value = 0.613495812322185
$ws.Range("E23").Value = 0.2214549941271144
$ws.Range("G23").Value = 0.002451927980081024
$ws.Range("I23").Value = 0.704893492283027
$ws.Range("J23").Value = 0.09669361577798696
$ws.Range("M23").Value = 0.6070641811384831
$ws.Range("N23").Value = 1.295503764027089
$ws.Range("O23").Value = 4.106845340936161
$ws.Range("B24").Value = 1.162888702100531
$ws.Range("C24").Value = 0.3160820911435565
$ws.Range("D24").Value = 0.6031305267781875
$ws.Range("E24").Value = 0.2196172115034756
$ws.Range("G24").Value = 0.00245687697355956
$ws.Range("I24").Value = 0.7094460264560141
$ws.Range("J24").Value = 0.09704940915504423
$ws.Range("M24").Value = 0.56169563579369
$ws.Range("N24").Value = 1.317886744104928
$ws.Range("O24").Value = 4.071539665233928
$ws.Range("B25").Value = 1.001494459348976
$ws.Range("C25").Value = 0.269804880993064
$ws.Range("D25").Value = 0.5930152468677221
$ws.Range("E25").Value = 0.2180168664246338
$ws.Range("G25").Value = 0.00246262026867938
$ws.Range("I25").Value = 0.7160182627004801
$ws.Range("J25").Value = 0.09759034889344065
$ws.Range("M25").Value = 0.5134907970860354
$ws.Range("N25").Value = 1.344043250370497
$ws.Range("O25").Value = 4.042977324913409
